$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new worksheet row at row 61 (shifts existing rows 61:69 down to 62:70)
$ws.Rows.Item(61).Insert()

# Grow the Excel Table (ListObject) to include the newly inserted row,
# keeping its range in sync with the new used range (A1:I70)
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:I70"))

# Fill in the new row's data for part #85 - DROK 5A buck converter mount
$ws.Cells.Item(61, 1).Value = 85
$ws.Cells.Item(61, 2).Value = "Electrical"
$ws.Cells.Item(61, 3).Value = "Mount"
$ws.Cells.Item(61, 4).Value = "N"
$ws.Cells.Item(61, 5).Value = "5A Buck Mount"
$ws.Cells.Item(61, 6).Value = "ABS"
$ws.Cells.Item(61, 7).Value = 1
$ws.Cells.Item(61, 9).Value = "85 - Electrical - Mount - 5A Buck Mount.stl"

# Move the active selection to match the post-edit workbook state
$ws.Range("I61").Select()
